# Regenerate orders with updated distance/size codes.
# The experiment's distance and stimulus-size labels changed:
#   D80 -> D86, D64 -> D69, D51 -> D55, S30 -> S31
# These substrings appear embedded throughout several text columns
# (Condition, Filename_Left, Filename_Right, Distance, Size), so walk
# the used range and rewrite any cell text containing them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ur = $ws.UsedRange
$rowCount = $ur.Rows.Count
$colCount = $ur.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($val -ne $null) {
            $sval = [string]$val
            $newval = $sval.Replace("D80", "D86").Replace("D64", "D69").Replace("D51", "D55").Replace("S30", "S31")
            if ($newval -ne $sval) {
                $cell.Value = $newval
            }
        }
    }
}
